# Auto-generated Excel COM-interop script to apply cryptos.xlsx updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '26.764.88'
$ws.Range("E2").Value = '  +0.88%  '
Set-TextCell $ws.Range("D3") '1.648.18'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("E4").Value = '  +0.36%  '
Set-TextCell $ws.Range("D5") '216.51'
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("E9").Value = '  +0.55%  '
Set-TextCell $ws.Range("D10") '19.25'
$ws.Range("E10").Value = '  +2.34%  '
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("E12").Value = '  +1.26%  '
Set-TextCell $ws.Range("B13") 'Polkadot'
Set-TextCell $ws.Range("C13") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws.Range("D13") '4.19'
$ws.Range("E13").Value = '  +1.24%  '
Set-TextCell $ws.Range("B14") 'WrappedEther'
Set-TextCell $ws.Range("C14") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws.Range("D14") '1.634.31'
$ws.Range("E14").Value = '  +0.73%  '
Set-TextCell $ws.Range("D15") '0.532'
$ws.Range("E15").Value = '  +1.62%  '
$ws.Range("E16").Value = '  +0.34%  '
Set-TextCell $ws.Range("D17") '26.760.45'
$ws.Range("E17").Value = '  +0.85%  '
$ws.Range("E18").Value = '  +0.35%  '
Set-TextCell $ws.Range("D19") '217.66'
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("E21").Value = '  +1.77%  '
$ws.Range("E22").Value = '  +13.61%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  +1.69%  '
Set-TextCell $ws.Range("D25") '146.76'
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("E27").Value = '  -0.15%  '
Set-TextCell $ws.Range("D28") '7.17'
$ws.Range("E28").Value = '  +3.93%  '
Set-TextCell $ws.Range("D29") '15.76'
$ws.Range("E29").Value = '  +1.17%  '
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("E33").Value = '  +1.46%  '
Set-TextCell $ws.Range("D34") '1.281.10'
$ws.Range("E34").Value = '  +3.35%  '
$ws.Range("E35").Value = '  +2.86%  '
$ws.Range("E36").Value = '  +2.76%  '
$ws.Range("E37").Value = '  +1.97%  '
Set-TextCell $ws.Range("D38") '0.538'
$ws.Range("E38").Value = '  +5.54%  '
Set-TextCell $ws.Range("D39") '0.830'
$ws.Range("E39").Value = '  +4.43%  '
$ws.Range("E40").Value = '  +0.36%  '
Set-TextCell $ws.Range("D41") '0.816'
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("E43").Value = '  +2.06%  '
Set-TextCell $ws.Range("D44") '1.788.57'
$ws.Range("E44").Value = '  +1.35%  '
Set-TextCell $ws.Range("D45") '91.98'
$ws.Range("E45").Value = '  -1.25%  '
Set-TextCell $ws.Range("D46") '59.75'
$ws.Range("E46").Value = '  +8.92%  '
Set-TextCell $ws.Range("D47") '1.60'
$ws.Range("E47").Value = '  +1.13%  '
$ws.Range("E48").Value = '  +0.19%  '
Set-TextCell $ws.Range("D49") '0.0516'
$ws.Range("E49").Value = '  +1.14%  '
Set-TextCell $ws.Range("D50") '7.78'
$ws.Range("E50").Value = '  +3.82%  '
Set-TextCell $ws.Range("D51") '0.0973'
$ws.Range("E51").Value = '  +1.57%  '
